$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 111702506
$ws.Range("B8").Value = 90687
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 5964
$ws.Range("F8").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G8").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H8").Value = '(L.:Fr.) P.Karst.'
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("P8").Value = 'Kyrkogården, Nrk'
$ws.Range("Q8").Value = 517093.6249861007
$ws.Range("R8").Value = 6574959.965416327

# Row 9
$ws.Range("A9").Value = 111702400
$ws.Range("B9").Value = 90687
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 5964
$ws.Range("F9").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G9").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H9").Value = '(L.:Fr.) P.Karst.'
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("P9").Value = 'Kyrkogården (Kyrkogården), Nrk'
$ws.Range("Q9").Value = 517073.2951468225
$ws.Range("R9").Value = 6574931.795150192

# Row 10
$ws.Range("A10").Value = 111702393
$ws.Range("B10").Value = 89183
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 3215
$ws.Range("F10").Value = 'Rödgul trumpetsvamp'
$ws.Range("G10").Value = 'Craterellus lutescens'
$ws.Range("H10").Value = '(Fr.) Fr.'
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("P10").Value = 'Kyrkogården (Kyrkogården), Nrk'
$ws.Range("Q10").Value = 517070.2129045375
$ws.Range("R10").Value = 6574934.844418272

# Row 11
$ws.Range("A11").Value = 111702486
$ws.Range("B11").Value = 90678
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 4366
$ws.Range("F11").Value = 'Skarp dropptaggsvamp'
$ws.Range("G11").Value = 'Hydnellum peckii'
$ws.Range("H11").Value = 'Banker'
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("P11").Value = 'Kyrkogården (Kyrkogården), Nrk'
$ws.Range("Q11").Value = 517080.8398438052
$ws.Range("R11").Value = 6574959.907818918

# Row 12
$ws.Range("A12").Value = 111702420
$ws.Range("B12").Value = 90709
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 5448
$ws.Range("F12").Value = 'Svartvit taggsvamp'
$ws.Range("G12").Value = 'Phellodon connatus'
$ws.Range("H12").Value = '(Schultz) nom.prov'
$ws.Range("I12").Value = "'1"
$ws.Range("J12").Value = 'fruktkroppar'
$ws.Range("P12").Value = 'Kyrkogården (Kyrkogården), Nrk'
$ws.Range("Q12").Value = 517086.1792710476
$ws.Range("R12").Value = 6574909.900584662
